$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "PSP-ID" values (e.g. 1.2.1, 1.1.1, ...) in column A of the work-package
# table (rows 26-41) are removed, leaving only the descriptive text in column B.
$ws.Range("A26:A41").ClearContents()

# Restore the selection/scroll state that the workbook was left in.
$ws.Range("I29:I30").Select()
